# B6-PowerPoint.pptx edit
#
# 1) The three data tables (one each on slides 14, 15 and 16) get a new
#    built-in table style applied (swapped from the custom "Table_0" style
#    to the built-in style {A9B8BD27-C58B-480D-946A-949D18EF9607}).
# 2) The presentation's theme colours are reset from the "Integral" /
#    "Red Violet" palette back to the stock Office palette (this is the
#    functional, COM-addressable part of the theme1.xml <-> theme2.xml
#    swap seen in the OOXML: the slide master's live theme -- reachable
#    through Slide.ThemeColorScheme -- is the only theme part PowerPoint's
#    object model exposes for editing; the Notes Master's theme part has
#    no COM surface, same as in real PowerPoint).

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------------
$newTableStyleId = "{A9B8BD27-C58B-480D-946A-949D18EF9607}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Restore the stock Office theme colours -----------------------------
function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    (RGB 0x00 0x00 0x00),  # 1  dk1
    (RGB 0xFF 0xFF 0xFF),  # 2  lt1
    (RGB 0x44 0x54 0x6A),  # 3  dk2
    (RGB 0xE7 0xE6 0xE6),  # 4  lt2
    (RGB 0x5B 0x9B 0xD5),  # 5  accent1
    (RGB 0xED 0x7D 0x31),  # 6  accent2
    (RGB 0xA5 0xA5 0xA5),  # 7  accent3
    (RGB 0xFF 0xC0 0x00),  # 8  accent4
    (RGB 0x44 0x72 0xC4),  # 9  accent5
    (RGB 0x70 0xAD 0x47),  # 10 accent6
    (RGB 0x05 0x63 0xC1),  # 11 hlink
    (RGB 0x95 0x4F 0x72)   # 12 folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
